$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows at the top of the price-history block (rows 72-74),
# pushing the existing rows 72:166 down to 75:169.
$ws.Rows("72:74").Insert()

# ---- New row 72: Extra quality entry for date 44966 ----
$ws.Range("A72").Value = 11
$ws.Range("B72").Value = "Vega Monumental Concepción"
$ws.Range("C72").Value = "Bíobío"
$ws.Range("D72").Value = 44966
$ws.Range("E72").Value = 8
$ws.Range("F72").Value = 100112028
$ws.Range("G72").Value = "Sandia"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Extra"
$ws.Range("J72").Value = 500
$ws.Range("K72").Value = 2800
$ws.Range("L72").Value = 2800
$ws.Range("M72").Value = 2800
$ws.Range("N72").Value = "$/unidad"
$ws.Range("O72").Value = "Región de O'Higgins"
$ws.Range("P72").Value = 2800
$ws.Range("Q72").Value = 1
$ws.Range("R72").Value = "Hortaliza"

# ---- New row 73: Primera quality entry for date 44966 ----
$ws.Range("A73").Value = 11
$ws.Range("B73").Value = "Vega Monumental Concepción"
$ws.Range("C73").Value = "Bíobío"
$ws.Range("D73").Value = 44966
$ws.Range("E73").Value = 8
$ws.Range("F73").Value = 100112028
$ws.Range("G73").Value = "Sandia"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 500
$ws.Range("K73").Value = 2200
$ws.Range("L73").Value = 2200
$ws.Range("M73").Value = 2200
$ws.Range("N73").Value = "$/unidad"
$ws.Range("O73").Value = "Región de O'Higgins"
$ws.Range("P73").Value = 2200
$ws.Range("Q73").Value = 1
$ws.Range("R73").Value = "Hortaliza"

# ---- New row 74: Segunda quality entry for date 44966 ----
$ws.Range("A74").Value = 11
$ws.Range("B74").Value = "Vega Monumental Concepción"
$ws.Range("C74").Value = "Bíobío"
$ws.Range("D74").Value = 44966
$ws.Range("E74").Value = 8
$ws.Range("F74").Value = 100112028
$ws.Range("G74").Value = "Sandia"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Segunda"
$ws.Range("J74").Value = 500
$ws.Range("K74").Value = 1800
$ws.Range("L74").Value = 1800
$ws.Range("M74").Value = 1800
$ws.Range("N74").Value = "$/unidad"
$ws.Range("O74").Value = "Región de O'Higgins"
$ws.Range("P74").Value = 1800
$ws.Range("Q74").Value = 1
$ws.Range("R74").Value = "Hortaliza"
